$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: role changes from "tutor" to "d" ---
$ws.Range("B7").Value = "d"

# --- Duplicate row 7's formatting/values down into the six new "paint" rows
#     (future rows 8,9,10,12,13,14 -- row 11 is a distinct, separately typed row
#     inserted in between) ---
$ws.Range("A7:G7").Copy($ws.Range("A8:G8"))
$ws.Range("A7:G7").Copy($ws.Range("A9:G9"))
$ws.Range("A7:G7").Copy($ws.Range("A10:G10"))
$ws.Range("A7:G7").Copy($ws.Range("A12:G12"))
$ws.Range("A7:G7").Copy($ws.Range("A13:G13"))
$ws.Range("A7:G7").Copy($ws.Range("A14:G14"))

# --- Fill in the per-row username / class_name values (row-number order so the
#     shared-string table is built up in the same order as the source edit) ---
$ws.Range("D8").Value = "tutor54@nkt.com"
$ws.Range("F8").Value = "paint sess 1"

$ws.Range("D9").Value = "tutor54@nkt.com"
$ws.Range("F9").Value = "paint fix 1"

$ws.Range("D10").Value = "tutor54@nkt.com"
$ws.Range("F10").Value = "paint var 1"

$ws.Range("D12").Value = "tutor54@nkt.com"
$ws.Range("F12").Value = "paint fix 1 ind"

$ws.Range("D13").Value = "tutor54@nkt.com"
$ws.Range("F13").Value = "paint var 1 ind"

$ws.Range("D14").Value = "tutor54@nkt.com"
# F14 keeps the value copied from F7 ("PHP session 4") -- its text is renamed
# to "Java session 8" below via the shared string that F7 also points to.

# --- New row 11: a distinct, independently-entered row (student message test) ---
$ws.Range("A11").Value = "signin"
$ws.Range("B11").Value = "tutor"
$ws.Range("C11").Value = "n"
$ws.Range("C11").HorizontalAlignment = -4131
$ws.Range("D11").Value = "srinivasesaivanan6324@gmail.com"
$ws.Range("E11").Value = "Test@1234"
$ws.Range("E11").Style = "Hyperlink"
$ws.Range("F11").Value = "clarinet session 4 multi"
$ws.Range("G11").Value = 5

# --- Rename the shared "PHP session 4" class name (F7/F14) to "Java session 8" ---
$ws.Range("F7").Value = "Java session 8"

# --- Hyperlinks for the new rows (username -> mailto:, password -> mailto:) ---
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:tutor54@nkt.com")
$ws.Hyperlinks.Add($ws.Range("E8"), "mailto:Admin@123")
$ws.Hyperlinks.Add($ws.Range("D9"), "mailto:tutor54@nkt.com")
$ws.Hyperlinks.Add($ws.Range("E9"), "mailto:Admin@123")
$ws.Hyperlinks.Add($ws.Range("D10"), "mailto:tutor54@nkt.com")
$ws.Hyperlinks.Add($ws.Range("E10"), "mailto:Admin@123")
$ws.Hyperlinks.Add($ws.Range("E11"), "mailto:Test@1234")
$ws.Hyperlinks.Add($ws.Range("D12"), "mailto:tutor54@nkt.com")
$ws.Hyperlinks.Add($ws.Range("E12"), "mailto:Admin@123")
$ws.Hyperlinks.Add($ws.Range("D13"), "mailto:tutor54@nkt.com")
$ws.Hyperlinks.Add($ws.Range("E13"), "mailto:Admin@123")
$ws.Hyperlinks.Add($ws.Range("D14"), "mailto:tutor54@nkt.com")
$ws.Hyperlinks.Add($ws.Range("E14"), "mailto:Admin@123")

# --- Selection moves to F11, matching the author's last edit location ---
$ws.Range("F11").Select()
